# TRE-210-BE: add "Tổng điểm tích lũy" (accumulated points) column to the
# revenue export template. A new column is inserted before the existing
# "Tổng doanh thu" column (H), shifting the trailing totals columns one
# slot to the right (H..L -> I..M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H, pushing the old H:L columns to I:M.
$ws.Columns.Item(8).Insert()

# Match the width of the new column to its left neighbour (both end up in
# the 16.33-wide group that spans F:H after the insert).
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# New header text for the inserted column.
$ws.Range("H8").Value = "Tổng điểm tích lũy"

# Reflect the new active selection/scroll position recorded for the sheet.
$ws.Activate() | Out-Null
$ws.Range("G14").Select() | Out-Null
